$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$dates = @(
    "09-08-2021",
    "10-08-2021",
    "11-08-2021",
    "12-08-2021",
    "13-08-2021",
    "16-08-2021",
    "17-08-2021",
    "18-08-2021",
    "19-08-2021",
    "20-08-2021",
    "23-08-2021",
    "24-08-2021",
    "25-08-2021",
    "26-08-2021",
    "27-08-2021",
    "30-08-2021",
    "31-08-2021",
    "01-09-2021",
    "02-09-2021",
    "03-09-2021",
    "06-09-2021",
    "07-09-2021"
)

# Reference style taken from an existing data row, used to make sure the
# newly written date cells end up with the same (default) cell style as
# all the other data rows, rather than a date-formatted style that Excel
# would otherwise auto-apply because the text looks like a date.
$refStyle = $ws.Cells.Item(2, 1).Style

$startRow = 152
for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $startRow + $i
    $cellA = $ws.Cells.Item($row, 1)

    # Force the cell to text format before assignment so Excel stores the
    # value as the literal date string instead of silently converting it
    # into a date serial number.
    $cellA.NumberFormat = "@"
    $cellA.Value = $dates[$i]
    # Restore the default (unformatted) style used by the rest of the data
    # rows now that the text has been entered as a string.
    $cellA.Style = $refStyle

    $ws.Cells.Item($row, 2).Value = 3.25
}
